# Auto-generated Excel COM-interop script
# Applies numeric cell updates (market-data refresh) across all 8 sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2733472.2
$ws.Range("J17").Value = 2733472.2
$ws.Range("L17").Value = 8200416.600000001
$ws.Range("N17").Value = -8200752.600000001
$ws.Range("H32").Value = 1195.2727
$ws.Range("J32").Value = 1285.4286
$ws.Range("L32").Value = 1285.4286
$ws.Range("N32").Value = -1937.4286
$ws.Range("H40").Value = 1888.7778
$ws.Range("J40").Value = 1571.2858
$ws.Range("L40").Value = 1571.2858
$ws.Range("N40").Value = -1921.2858
$ws.Range("H74").Value = 3648.4167
$ws.Range("I74").Value = 2878.1
$ws.Range("K74").Value = 2878.1
$ws.Range("M74").Value = -1942.1
$ws.Range("H77").Value = 3648.4167
$ws.Range("I77").Value = 2878.1
$ws.Range("K77").Value = 14390.5
$ws.Range("M77").Value = -9710.5
$ws.Range("H82").Value = 1316.8
$ws.Range("I82").Value = 1316.8
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 3950.4
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -3544.4
$ws.Range("N82").Value = ""
$ws.Range("H85").Value = 1316.8
$ws.Range("I85").Value = 1316.8
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 3950.4
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -2546.4
$ws.Range("N85").Value = ""
$ws.Range("H86").Value = 5570.05
$ws.Range("I86").Value = 5428.273
$ws.Range("J86").Value = 5743.3335
$ws.Range("K86").Value = 5428.273
$ws.Range("L86").Value = 5743.3335
$ws.Range("M86").Value = -4305.273
$ws.Range("N86").Value = -7989.3335
$ws.Range("H88").Value = 4152.4375
$ws.Range("I88").Value = 5174.8
$ws.Range("J88").Value = 3687.7273
$ws.Range("K88").Value = 5174.8
$ws.Range("L88").Value = 3687.7273
$ws.Range("M88").Value = -4768.8
$ws.Range("N88").Value = -4499.7273
$ws.Range("H89").Value = 5570.05
$ws.Range("I89").Value = 5428.273
$ws.Range("J89").Value = 5743.3335
$ws.Range("K89").Value = 27141.365
$ws.Range("L89").Value = 28716.6675
$ws.Range("M89").Value = -21525.365
$ws.Range("N89").Value = -39948.6675
$ws.Range("H91").Value = 4152.4375
$ws.Range("I91").Value = 5174.8
$ws.Range("J91").Value = 3687.7273
$ws.Range("K91").Value = 5174.8
$ws.Range("L91").Value = 3687.7273
$ws.Range("M91").Value = -3770.8
$ws.Range("N91").Value = -6495.7273
$ws.Range("H96").Value = 1260.1482
$ws.Range("I96").Value = 1112.9
$ws.Range("J96").Value = 1680.8572
$ws.Range("K96").Value = 3338.7
$ws.Range("L96").Value = 5042.571599999999
$ws.Range("M96").Value = -1965.7
$ws.Range("N96").Value = -7788.571599999999
$ws.Range("H99").Value = 1299.1428
$ws.Range("I99").Value = 1349
$ws.Range("J99").Value = 1000
$ws.Range("K99").Value = 4047
$ws.Range("L99").Value = 3000
$ws.Range("M99").Value = -2549
$ws.Range("N99").Value = -5996
$ws.Range("H100").Value = 2402.3225
$ws.Range("I100").Value = 1930
$ws.Range("K100").Value = 1930
$ws.Range("M100").Value = -1389
$ws.Range("H113").Value = 14907.714
$ws.Range("I113").Value = 41100
$ws.Range("K113").Value = 41100
$ws.Range("M113").Value = -37846
$ws.Range("H116").Value = 9248.333000000001
$ws.Range("J116").Value = 9872.5
$ws.Range("L116").Value = 9872.5
$ws.Range("N116").Value = -16756.5
$ws.Range("H121").Value = 3900
$ws.Range("J121").Value = 3900
$ws.Range("L121").Value = 11700
$ws.Range("N121").Value = -15194
$ws.Range("H135").Value = 2415
$ws.Range("I135").Value = 2415
$ws.Range("K135").Value = 21735
$ws.Range("M135").Value = -19200
$ws.Range("H138").Value = 2223.957
$ws.Range("I138").Value = 1561.7916
$ws.Range("J138").Value = 2454.2754
$ws.Range("K138").Value = 4685.3748
$ws.Range("L138").Value = 7362.8262
$ws.Range("M138").Value = 454.6252000000004
$ws.Range("N138").Value = -17642.8262
$ws.Range("H141").Value = 4619.778
$ws.Range("I141").Value = 5878.8335
$ws.Range("J141").Value = 2101.6667
$ws.Range("K141").Value = 17636.5005
$ws.Range("L141").Value = 6305.000100000001
$ws.Range("M141").Value = -12456.5005
$ws.Range("N141").Value = -16665.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 8480.200000000001
$ws.Range("I8").Value = 200
$ws.Range("J8").Value = 10550.25
$ws.Range("K8").Value = 200
$ws.Range("L8").Value = 10550.25
$ws.Range("M8").Value = -56
$ws.Range("N8").Value = -10838.25
$ws.Range("H10").Value = 135.66667
$ws.Range("I10").Value = 1
$ws.Range("J10").Value = 405
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 405
$ws.Range("M10").Value = 169
$ws.Range("N10").Value = -745
$ws.Range("H32").Value = 3694.54
$ws.Range("I32").Value = 2379.814
$ws.Range("K32").Value = 2379.814
$ws.Range("M32").Value = -2092.814
$ws.Range("H43").Value = 32999
$ws.Range("I43").Value = 31996
$ws.Range("J43").Value = 33249.75
$ws.Range("K43").Value = 31996
$ws.Range("L43").Value = 33249.75
$ws.Range("M43").Value = -31683
$ws.Range("N43").Value = -33875.75
$ws.Range("H45").Value = 2927.818
$ws.Range("I45").Value = 3058.125
$ws.Range("J45").Value = 2580.3333
$ws.Range("K45").Value = 3058.125
$ws.Range("L45").Value = 2580.3333
$ws.Range("M45").Value = -2681.125
$ws.Range("N45").Value = -3334.3333
$ws.Range("H61").Value = 68539.25999999999
$ws.Range("I61").Value = 1637.8182
$ws.Range("K61").Value = 1637.8182
$ws.Range("M61").Value = -1425.8182
$ws.Range("H74").Value = 10991.791
$ws.Range("I74").Value = 1665.4054
$ws.Range("J74").Value = 68504.5
$ws.Range("K74").Value = 1665.4054
$ws.Range("L74").Value = 68504.5
$ws.Range("M74").Value = -791.4054000000001
$ws.Range("N74").Value = -70252.5
$ws.Range("H77").Value = 10991.791
$ws.Range("I77").Value = 1665.4054
$ws.Range("J77").Value = 68504.5
$ws.Range("K77").Value = 8327.027
$ws.Range("L77").Value = 342522.5
$ws.Range("M77").Value = -3959.027
$ws.Range("N77").Value = -351258.5
$ws.Range("H82").Value = 30000
$ws.Range("J82").Value = 30000
$ws.Range("L82").Value = 30000
$ws.Range("N82").Value = -30722
$ws.Range("H85").Value = 30000
$ws.Range("J85").Value = 30000
$ws.Range("L85").Value = 30000
$ws.Range("N85").Value = -32496
$ws.Range("H88").Value = 1706.4166
$ws.Range("I88").Value = 1185.3334
$ws.Range("K88").Value = 1185.3334
$ws.Range("M88").Value = -779.3334
$ws.Range("H91").Value = 1706.4166
$ws.Range("I91").Value = 1185.3334
$ws.Range("K91").Value = 1185.3334
$ws.Range("M91").Value = 218.6666
$ws.Range("H97").Value = 614.9355
$ws.Range("I97").Value = 553.76
$ws.Range("K97").Value = 553.76
$ws.Range("M97").Value = -57.75999999999999
$ws.Range("H102").Value = 4753.7856
$ws.Range("I102").Value = 2755.8
$ws.Range("J102").Value = 9748.75
$ws.Range("K102").Value = 2755.8
$ws.Range("L102").Value = 9748.75
$ws.Range("M102").Value = -1133.8
$ws.Range("N102").Value = -12992.75
$ws.Range("H122").Value = 1093223.2
$ws.Range("I122").Value = 1595688.9
$ws.Range("J122").Value = 4547.6665
$ws.Range("K122").Value = 4787066.699999999
$ws.Range("L122").Value = 13642.9995
$ws.Range("M122").Value = -4784616.699999999
$ws.Range("N122").Value = -18542.9995
$ws.Range("H136").Value = 68539.25999999999
$ws.Range("I136").Value = 1637.8182
$ws.Range("K136").Value = 4913.4546
$ws.Range("M136").Value = -2363.4546

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 6683980
$ws.Range("I20").Value = 13893921
$ws.Range("K20").Value = 13893921
$ws.Range("M20").Value = -13893674
$ws.Range("H86").Value = 47621010
$ws.Range("I86").Value = 1731.5333
$ws.Range("J86").Value = 166669200
$ws.Range("K86").Value = 1731.5333
$ws.Range("L86").Value = 166669200
$ws.Range("M86").Value = -608.5333000000001
$ws.Range("N86").Value = -166671446
$ws.Range("H89").Value = 47621010
$ws.Range("I89").Value = 1731.5333
$ws.Range("J89").Value = 166669200
$ws.Range("K89").Value = 8657.666499999999
$ws.Range("L89").Value = 833346000
$ws.Range("M89").Value = -3041.666499999999
$ws.Range("N89").Value = -833357232
$ws.Range("H94").Value = 4753.2856
$ws.Range("I94").Value = 2127.4443
$ws.Range("J94").Value = 9479.799999999999
$ws.Range("K94").Value = 2127.4443
$ws.Range("L94").Value = 9479.799999999999
$ws.Range("M94").Value = -1676.4443
$ws.Range("N94").Value = -10381.8
$ws.Range("H99").Value = 21689.826
$ws.Range("J99").Value = 4624.5
$ws.Range("L99").Value = 4624.5
$ws.Range("N99").Value = -7620.5
$ws.Range("H134").Value = 34002.047
$ws.Range("J134").Value = 24187.4
$ws.Range("L134").Value = 72562.20000000001
$ws.Range("N134").Value = -77632.20000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 9751.4
$ws.Range("I16").Value = 7303.7
$ws.Range("K16").Value = 7303.7
$ws.Range("M16").Value = -7016.7
$ws.Range("H31").Value = 16212.637
$ws.Range("I31").Value = 1593.375
$ws.Range("J31").Value = 24566.5
$ws.Range("K31").Value = 1593.375
$ws.Range("L31").Value = 24566.5
$ws.Range("M31").Value = -1298.375
$ws.Range("N31").Value = -25156.5
$ws.Range("H34").Value = 16212.637
$ws.Range("I34").Value = 1593.375
$ws.Range("J34").Value = 24566.5
$ws.Range("K34").Value = 1593.375
$ws.Range("L34").Value = 24566.5
$ws.Range("M34").Value = -1391.375
$ws.Range("N34").Value = -24970.5
$ws.Range("H76").Value = 7850
$ws.Range("I76").Value = 7850
$ws.Range("K76").Value = 7850
$ws.Range("M76").Value = -7535
$ws.Range("H79").Value = 7850
$ws.Range("I79").Value = 7850
$ws.Range("K79").Value = 7850
$ws.Range("M79").Value = -6758
$ws.Range("H86").Value = 14043.75
$ws.Range("I86").Value = 16173.223
$ws.Range("K86").Value = 16173.223
$ws.Range("M86").Value = -15050.223
$ws.Range("H89").Value = 14043.75
$ws.Range("I89").Value = 16173.223
$ws.Range("K89").Value = 80866.11500000001
$ws.Range("M89").Value = -75250.11500000001
$ws.Range("H107").Value = 1424.2
$ws.Range("I107").Value = 1570.2727
$ws.Range("K107").Value = 1570.2727
$ws.Range("M107").Value = 349.7273
$ws.Range("H113").Value = 9751.4
$ws.Range("I113").Value = 7303.7
$ws.Range("K113").Value = 7303.7
$ws.Range("M113").Value = -5133.7
$ws.Range("H125").Value = 84784.75
$ws.Range("J125").Value = 84784.75
$ws.Range("L125").Value = 84784.75
$ws.Range("N125").Value = -89704.75
$ws.Range("H132").Value = 5930
$ws.Range("I132").Value = 4950
$ws.Range("K132").Value = 14850
$ws.Range("M132").Value = -12320
$ws.Range("H134").Value = 33339704
$ws.Range("I134").Value = 1590.4286
$ws.Range("K134").Value = 4771.2858
$ws.Range("M134").Value = -2236.2858

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 101.77778
$ws.Range("J2").Value = 43.692307
$ws.Range("L2").Value = 262.153842
$ws.Range("N2").Value = -488.153842
$ws.Range("H5").Value = 11150973
$ws.Range("J5").Value = 27875670
$ws.Range("L5").Value = 83627010
$ws.Range("N5").Value = -83627234
$ws.Range("H34").Value = 3362.4
$ws.Range("I34").Value = 379.83334
$ws.Range("J34").Value = 7836.25
$ws.Range("K34").Value = 1139.50002
$ws.Range("L34").Value = 23508.75
$ws.Range("M34").Value = -1055.50002
$ws.Range("N34").Value = -23676.75
$ws.Range("H51").Value = 335.33334
$ws.Range("I51").Value = 68
$ws.Range("J51").Value = 602.6667
$ws.Range("K51").Value = 204
$ws.Range("L51").Value = 1808.0001
$ws.Range("M51").Value = 256
$ws.Range("N51").Value = -2728.0001
$ws.Range("H61").Value = 352.5
$ws.Range("I61").Value = 383
$ws.Range("J61").Value = 200
$ws.Range("K61").Value = 1149
$ws.Range("L61").Value = 600
$ws.Range("M61").Value = -934
$ws.Range("N61").Value = -1030
$ws.Range("H92").Value = 345.58334
$ws.Range("I92").Value = 335.7143
$ws.Range("J92").Value = 359.4
$ws.Range("K92").Value = 1007.1429
$ws.Range("L92").Value = 1078.2
$ws.Range("M92").Value = 240.8571000000001
$ws.Range("N92").Value = -3574.2
$ws.Range("H99").Value = 1756.25
$ws.Range("I99").Value = 512.5
$ws.Range("J99").Value = 3000
$ws.Range("K99").Value = 1537.5
$ws.Range("L99").Value = 9000
$ws.Range("M99").Value = 708.5
$ws.Range("N99").Value = -13492
$ws.Range("H107").Value = 1079920.5
$ws.Range("I107").Value = 539.6923
$ws.Range("K107").Value = 1619.0769
$ws.Range("M107").Value = 300.9231
$ws.Range("H112").Value = 13307.692
$ws.Range("I112").Value = 5000
$ws.Range("J112").Value = 14000
$ws.Range("K112").Value = 15000
$ws.Range("L112").Value = 42000
$ws.Range("M112").Value = -13892
$ws.Range("N112").Value = -44216
$ws.Range("H129").Value = 5682946.5
$ws.Range("J129").Value = 9091809
$ws.Range("L129").Value = 27275427
$ws.Range("N129").Value = -27285427
$ws.Range("H131").Value = 1453.09
$ws.Range("J131").Value = 1453.09
$ws.Range("L131").Value = 4359.27
$ws.Range("N131").Value = -14439.27
$ws.Range("H134").Value = 5952.974
$ws.Range("I134").Value = 2405.8
$ws.Range("K134").Value = 7217.400000000001
$ws.Range("M134").Value = -2147.400000000001
$ws.Range("H135").Value = 11150973
$ws.Range("J135").Value = 27875670
$ws.Range("L135").Value = 250881030
$ws.Range("N135").Value = -250886100
$ws.Range("H137").Value = 3872.25
$ws.Range("I137").Value = 4333
$ws.Range("J137").Value = 3595.8
$ws.Range("K137").Value = 12999
$ws.Range("L137").Value = 10787.4
$ws.Range("M137").Value = -7899
$ws.Range("N137").Value = -20987.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 208.58824
$ws.Range("I2").Value = 215.83333
$ws.Range("J2").Value = 204.63637
$ws.Range("K2").Value = 215.83333
$ws.Range("L2").Value = 204.63637
$ws.Range("M2").Value = -102.83333
$ws.Range("N2").Value = -430.63637
$ws.Range("H7").Value = 275037.25
$ws.Range("I7").Value = 33383
$ws.Range("K7").Value = 33383
$ws.Range("M7").Value = -33271
$ws.Range("H8").Value = 275037.25
$ws.Range("I8").Value = 33383
$ws.Range("K8").Value = 33383
$ws.Range("M8").Value = -33244
$ws.Range("H18").Value = 22674.5
$ws.Range("I18").Value = 15349
$ws.Range("K18").Value = 15349
$ws.Range("M18").Value = -15056
$ws.Range("H70").Value = 11752
$ws.Range("I70").Value = 12499.5
$ws.Range("J70").Value = 11502.833
$ws.Range("K70").Value = 12499.5
$ws.Range("L70").Value = 11502.833
$ws.Range("M70").Value = -12229.5
$ws.Range("N70").Value = -12042.833
$ws.Range("H73").Value = 11752
$ws.Range("I73").Value = 12499.5
$ws.Range("J73").Value = 11502.833
$ws.Range("K73").Value = 12499.5
$ws.Range("L73").Value = 11502.833
$ws.Range("M73").Value = -11563.5
$ws.Range("N73").Value = -13374.833
$ws.Range("H80").Value = 13741.211
$ws.Range("I80").Value = 12115.5
$ws.Range("J80").Value = 16528.143
$ws.Range("K80").Value = 12115.5
$ws.Range("L80").Value = 16528.143
$ws.Range("M80").Value = -11117.5
$ws.Range("N80").Value = -18524.143
$ws.Range("H83").Value = 13741.211
$ws.Range("I83").Value = 12115.5
$ws.Range("J83").Value = 16528.143
$ws.Range("K83").Value = 60577.5
$ws.Range("L83").Value = 82640.715
$ws.Range("M83").Value = -55585.5
$ws.Range("N83").Value = -92624.715
$ws.Range("H131").Value = 50000
$ws.Range("J131").Value = 50000
$ws.Range("L131").Value = 50000
$ws.Range("N131").Value = -60080
$ws.Range("H132").Value = 3475.6924
$ws.Range("I132").Value = 3347.5
$ws.Range("K132").Value = 10042.5
$ws.Range("M132").Value = -7512.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 25000
$ws.Range("J3").Value = 25000
$ws.Range("L3").Value = 25000
$ws.Range("N3").Value = -25224
$ws.Range("H15").Value = 25000
$ws.Range("J15").Value = 25000
$ws.Range("L15").Value = 25000
$ws.Range("N15").Value = -25340
$ws.Range("H16").Value = 76925140
$ws.Range("I16").Value = 111113800
$ws.Range("J16").Value = 648.5
$ws.Range("K16").Value = 111113800
$ws.Range("L16").Value = 648.5
$ws.Range("M16").Value = -111113630
$ws.Range("N16").Value = -988.5
$ws.Range("H35").Value = 5906.8
$ws.Range("I35").Value = 2249.5
$ws.Range("K35").Value = 2249.5
$ws.Range("M35").Value = -1913.5
$ws.Range("H40").Value = 1686246.5
$ws.Range("I40").Value = 4568.636
$ws.Range("J40").Value = 4532163
$ws.Range("K40").Value = 4568.636
$ws.Range("L40").Value = 4532163
$ws.Range("M40").Value = -4432.636
$ws.Range("N40").Value = -4532435
$ws.Range("H46").Value = 2707.5833
$ws.Range("I46").Value = 1299
$ws.Range("J46").Value = 2835.6365
$ws.Range("K46").Value = 1299
$ws.Range("L46").Value = 2835.6365
$ws.Range("M46").Value = -1111
$ws.Range("N46").Value = -3211.6365
$ws.Range("H61").Value = 2635.1614
$ws.Range("I61").Value = 2026.5769
$ws.Range("K61").Value = 2026.5769
$ws.Range("M61").Value = -1824.5769
$ws.Range("H69").Value = 80000
$ws.Range("J69").Value = 80000
$ws.Range("L69").Value = 80000
$ws.Range("N69").Value = -81622
$ws.Range("H72").Value = 80000
$ws.Range("J72").Value = 80000
$ws.Range("L72").Value = 240000
$ws.Range("N72").Value = -248112
$ws.Range("H81").Value = 100000
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 100000
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 100000
$ws.Range("N81").Value = -101996
$ws.Range("M81").Value = ""
$ws.Range("H82").Value = 1928.3024
$ws.Range("I82").Value = 2156.52
$ws.Range("J82").Value = 1611.3334
$ws.Range("K82").Value = 2156.52
$ws.Range("L82").Value = 1611.3334
$ws.Range("M82").Value = -1795.52
$ws.Range("N82").Value = -2333.3334
$ws.Range("H84").Value = 100000
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 100000
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 300000
$ws.Range("N84").Value = -309984
$ws.Range("M84").Value = ""
$ws.Range("H85").Value = 1928.3024
$ws.Range("I85").Value = 2156.52
$ws.Range("J85").Value = 1611.3334
$ws.Range("K85").Value = 2156.52
$ws.Range("L85").Value = 1611.3334
$ws.Range("M85").Value = -908.52
$ws.Range("N85").Value = -4107.3334
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").Value = ""
$ws.Range("H104").Value = 13841
$ws.Range("J104").Value = 13841
$ws.Range("L104").Value = 13841
$ws.Range("N104").Value = -20829
$ws.Range("H113").Value = 2635.1614
$ws.Range("I113").Value = 2026.5769
$ws.Range("K113").Value = 2026.5769
$ws.Range("M113").Value = 143.4231
$ws.Range("H122").Value = 28727628
$ws.Range("J122").Value = 3131136.8
$ws.Range("L122").Value = 9393410.399999999
$ws.Range("N122").Value = -9398310.399999999
$ws.Range("H136").Value = 9543.736000000001
$ws.Range("I136").Value = 6990.4614
$ws.Range("K136").Value = 20971.3842
$ws.Range("M136").Value = -18421.3842

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 79999
$ws.Range("J16").Value = 79999
$ws.Range("L16").Value = 79999
$ws.Range("N16").Value = -80583
$ws.Range("H76").Value = 46778
$ws.Range("J76").Value = 46778
$ws.Range("L76").Value = 46778
$ws.Range("N76").Value = -47408
$ws.Range("H79").Value = 46778
$ws.Range("J79").Value = 46778
$ws.Range("L79").Value = 46778
$ws.Range("N79").Value = -48962
$ws.Range("H81").Value = 3496.2144
$ws.Range("I81").Value = 3672.8462
$ws.Range("K81").Value = 7345.6924
$ws.Range("M81").Value = -6284.6924
$ws.Range("H84").Value = 3496.2144
$ws.Range("I84").Value = 3672.8462
$ws.Range("K84").Value = 36728.462
$ws.Range("M84").Value = -31424.462
$ws.Range("H100").Value = 378.72726
$ws.Range("J100").Value = 451.57144
$ws.Range("L100").Value = 903.14288
$ws.Range("N100").Value = -1985.14288
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").Value = ""
$ws.Range("H136").Value = 14904.546
$ws.Range("I136").Value = 2377.0833
$ws.Range("K136").Value = 7131.249899999999
$ws.Range("M136").Value = -4581.249899999999
